$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the existing row 679 (shifts old 679..720 down to 681..722,
# and Excel auto-extends the used range accordingly).
$ws.Rows("679:680").Insert()

function Set-DateText($row, $dateText, $dow, $hour) {
    $a = $ws.Cells.Item($row, 1)
    $a.NumberFormat = "@"
    $a.Value = $dateText
    $a.ClearFormats()
    $ws.Cells.Item($row, 2).Value = $dow
    $ws.Cells.Item($row, 3).Value = $hour
    $ws.Cells.Item($row, 4).Value = 201
}

Set-DateText 679 "2026/01/20" "火" 22
Set-DateText 680 "2026/01/21" "水" 2
